# Updates the "cryptos" price/volume snapshot (Price column D, Volume(1h)
# column E), matching the GitHub Actions refresh commit. A couple of cells
# in column D (e.g. "1.000") look like plain numbers to Excel's parser, so
# those are written with a leading "'" to force text, then the cell style
# is reset back to "Normal" so no stray number-format/style is left behind
# on the cell (keeping it identical to the untouched text cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.429.47'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").Value = '1.639.68'
$ws.Range("E3").Value = '  +2.18%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").Value = "'305.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.49%  '
$ws.Range("D7").Value = "'0.3727"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.41%  '
$ws.Range("D8").Value = "'52.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.73%  '
$ws.Range("D9").Value = "'0.3622"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.31%  '
$ws.Range("D10").Value = "'1.252"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.63%  '
$ws.Range("D11").Value = "'0.08119"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").Value = "'1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("E13").Value = '  -0.11%  '
$ws.Range("D14").Value = "'6.598"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("E15").Value = '  +1.65%  '
$ws.Range("D16").Value = "'7.285"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.84%  '
$ws.Range("D17").Value = '1.630.05'
$ws.Range("E17").Value = '  +1.70%  '
$ws.Range("D18").Value = "'94.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.46%  '
$ws.Range("D19").Value = "'0.06873"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.39%  '
$ws.Range("D20").Value = "'18.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D22").Value = "'0.9994"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '23.440.52'
$ws.Range("E23").Value = '  +0.89%  '
$ws.Range("E24").Value = '  -1.94%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = "'2.410"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.11%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = "'3.019"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").Value = "'21.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").Value = "'151.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.08%  '
$ws.Range("D29").Value = "'5.289"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.80%  '
$ws.Range("D30").Value = "'135.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.03%  '
$ws.Range("D31").Value = "'2.284"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.51%  '
$ws.Range("D32").Value = '1.811.49'
$ws.Range("D33").Value = "'6.752"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.26%  '
$ws.Range("D34").Value = "'0.9543"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.30%  '
$ws.Range("D35").Value = "'0.02834"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.97%  '
$ws.Range("D36").Value = "'10.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.69%  '
$ws.Range("D37").Value = "'0.2520"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("D38").Value = "'0.07215"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.90%  '
$ws.Range("D39").Value = "'0.08783"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.27%  '
$ws.Range("D40").Value = "'6.048"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.81%  '
$ws.Range("D41").Value = "'1.375"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.36%  '
$ws.Range("D42").Value = "'0.7043"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("D43").Value = "'12.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.98%  '
$ws.Range("D44").Value = "'16.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.00%  '
$ws.Range("D45").Value = "'0.6510"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("D46").Value = "'2.323"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").Value = "'0.9996"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("D48").Value = "'4.007"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.48%  '
$ws.Range("D49").Value = "'0.07973"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("D50").Value = "'128.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.08%  '
$ws.Range("D51").Value = "'1.197"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.77%  '
